$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.826.07'
$ws.Range('E2').Value = '  -1.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.633.46'
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.53'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5020'
$ws.Range('E6').Value = '  -1.42%  '
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2565'
$ws.Range('E8').Value = '  -0.62%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06384'
$ws.Range('E9').Value = '  -0.47%  '
$ws.Range('E10').Value = '  -1.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07689'
$ws.Range('E11').Value = '  -1.50%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.637.46'
$ws.Range('E12').Value = '  -1.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.242'
$ws.Range('E13').Value = '  -0.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.858.00'
$ws.Range('E15').Value = '  -1.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅7913'
$ws.Range('E16').Value = '  -1.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.52'
$ws.Range('E17').Value = '  -0.64%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.842.81'
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.002'
$ws.Range('E19').Value = '  -0.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '201.77'
$ws.Range('E20').Value = '  -3.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.326'
$ws.Range('E21').Value = '  -1.80%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.916'
$ws.Range('E22').Value = '  -1.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.957'
$ws.Range('E23').Value = '  -0.76%  '
$ws.Range('E24').Value = '  -0.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.933'
$ws.Range('E25').Value = '  +11.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.68'
$ws.Range('E26').Value = '  -1.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1139'
$ws.Range('E27').Value = '  -2.72%  '
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('E29').Value = '  -3.78%  '
$ws.Range('E30').Value = '  -0.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04998'
$ws.Range('E31').Value = '  -2.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.262'
$ws.Range('E32').Value = '  -2.55%  '
$ws.Range('E33').Value = '  -1.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.540'
$ws.Range('E34').Value = '  -1.72%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.366'
$ws.Range('E35').Value = '  -0.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.172.24'
$ws.Range('E36').Value = '  +1.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.8907'
$ws.Range('E37').Value = '  -4.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.614'
$ws.Range('E38').Value = '  -4.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5584'
$ws.Range('E39').Value = '  -1.90%  '
$ws.Range('E40').Value = '  -1.83%  '
$ws.Range('E41').Value = '  -0.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.684'
$ws.Range('E42').Value = '  +0.84%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8054'
$ws.Range('E43').Value = '  -3.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.38'
$ws.Range('E44').Value = '  -0.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.770.01'
$ws.Range('E45').Value = '  -1.23%  '
$ws.Range('E46').Value = '  -1.14%  '
$ws.Range('E47').Value = '  -0.76%  '
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range('E49').Value = '  -1.66%  '
$ws.Range('E50').Value = '  +0.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.485'
$ws.Range('E51').Value = '  -4.88%  '
